# "un tas de nouveau classement depuis autre pc"
# Fill in the "Inscriptions via Dossardeur" placeholder in column E with the
# real per-event registration link keys for a batch of newly-known events.
# New shared strings must be created in this exact order so they land at
# sharedStrings indices 211-217 (morschwiller_vtt, heimsbrunn_cx, wittenheim,
# blaesheim, technochape, boron, Frotey_clm) to match the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calendrier")

$ws.Range("E55").Value = "morschwiller_vtt"
$ws.Range("E56").Value = "heimsbrunn_cx"
$ws.Range("E47").Value = "wittenheim"
$ws.Range("E52").Value = "blaesheim"
$ws.Range("E48").Value = "technochape"
$ws.Range("E43").Value = "boron"
$ws.Range("E45").Value = "Frotey_clm"

# Match the author's final cursor position / view (best-effort; scroll
# position itself is not persisted by this runtime but the selection is).
$ws.Range("C48").Select()
